# Refresh the cryptos list per the "Updated cryptos list ... with GitHub Actions" commit.
# Rewrites the changed Coin/Link/Price/Volume(1h) cells (rows 2-51) to their new values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.700.81'
$ws.Range('E2').Value = '  +0.22%  '
$ws.Range('D3').Value = '3.737.45'
$ws.Range('E3').Value = '  +7.07%  '
$ws.Range('E4').Value = '  +0.12%  '
$ws.Range('D5').Value = '''613.01'
$ws.Range('E5').Value = '  +4.30%  '
$ws.Range('D6').Value = '''177.61'
$ws.Range('E6').Value = '  -4.12%  '
$ws.Range('D7').Value = '3.735.68'
$ws.Range('E7').Value = '  +6.18%  '
$ws.Range('E8').Value = '  +0.37%  '
$ws.Range('E9').Value = '  +1.55%  '
$ws.Range('E10').Value = '  +5.68%  '
$ws.Range('E11').Value = '  -3.65%  '
$ws.Range('D12').Value = '''0.498'
$ws.Range('E12').Value = '  +1.87%  '
$ws.Range('D13').Value = '''40.91'
$ws.Range('E13').Value = '  +5.89%  '
$ws.Range('E14').Value = '  +1.73%  '
$ws.Range('D15').Value = '4.361.59'
$ws.Range('E15').Value = '  +7.89%  '
$ws.Range('D16').Value = '3.737.51'
$ws.Range('E16').Value = '  +8.27%  '
$ws.Range('D17').Value = '69.746.87'
$ws.Range('E17').Value = '  -0.03%  '
$ws.Range('E18').Value = '  +0.50%  '
$ws.Range('D19').Value = '''7.60'
$ws.Range('E19').Value = '  +2.58%  '
$ws.Range('D20').Value = '''514.84'
$ws.Range('E20').Value = '  +2.44%  '
$ws.Range('D21').Value = '''16.74'
$ws.Range('E21').Value = '  -1.47%  '
$ws.Range('D22').Value = '''9.58'
$ws.Range('E22').Value = '  +8.40%  '
$ws.Range('E23').Value = '  -0.27%  '
$ws.Range('D24').Value = '''88.12'
$ws.Range('E24').Value = '  +1.90%  '
$ws.Range('D25').Value = '''2.51'
$ws.Range('E25').Value = '  +5.69%  '
$ws.Range('D26').Value = '''13.39'
$ws.Range('E26').Value = '  +0.71%  '
$ws.Range('D27').Value = '''11.08'
$ws.Range('E27').Value = '  +3.34%  '
$ws.Range('E28').Value = '  -0.04%  '
$ws.Range('E29').Value = '  +18.71%  '
$ws.Range('E30').Value = '  +0.28%  '
$ws.Range('B31').Value = 'PancakeSwap'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D31').Value = '''2.84'
$ws.Range('E31').Value = '  +4.86%  '
$ws.Range('B32').Value = 'NEARProtocol'
$ws.Range('C32').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D32').Value = '''7.85'
$ws.Range('E32').Value = '  -3.54%  '
$ws.Range('E33').Value = '  +2.04%  '
$ws.Range('E34').Value = '  -0.82%  '
$ws.Range('E35').Value = '  -0.10%  '
$ws.Range('D36').Value = '''6.22'
$ws.Range('E36').Value = '  +1.54%  '
$ws.Range('E37').Value = '  +2.54%  '
$ws.Range('E38').Value = '  +3.28%  '
$ws.Range('E39').Value = '  +3.56%  '
$ws.Range('E40').Value = '  +4.24%  '
$ws.Range('D41').Value = '''51.32'
$ws.Range('E41').Value = '  +2.38%  '
$ws.Range('B42').Value = 'Arweave'
$ws.Range('C42').Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range('D42').Value = '''44.39'
$ws.Range('E42').Value = '  -7.38%  '
$ws.Range('B43').Value = 'Cosmos'
$ws.Range('C43').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D43').Value = '''8.83'
$ws.Range('E43').Value = '  +1.33%  '
$ws.Range('D44').Value = '''424.38'
$ws.Range('E44').Value = '  +3.97%  '
$ws.Range('D45').Value = '3.093.46'
$ws.Range('E45').Value = '  +4.28%  '
$ws.Range('E46').Value = '  -3.94%  '
$ws.Range('D47').Value = '''0.0364'
$ws.Range('E47').Value = '  +1.26%  '
$ws.Range('D48').Value = '''27.87'
$ws.Range('E48').Value = '  -1.09%  '
$ws.Range('D49').Value = '''2.53'
$ws.Range('E49').Value = '  +3.31%  '
$ws.Range('D51').Value = '''135.35'
$ws.Range('E51').Value = '  +0.10%  '
